$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: now the "Appointment Fee" charge row ---
$ws.Cells.Item(2, 2).Value = "APPT FEE"
$ws.Cells.Item(2, 3).Value = "Appointment Fee"
$ws.Cells.Item(2, 4).Value = "Missed or Uncheduled Delivery Appointment"
$ws.Range("E2:I2").ClearContents()

# --- Row 3: now the "OverFlow:InboundPallets" charge row ---
$ws.Cells.Item(3, 2).Value = "OverFlow:InboundPallets"
$ws.Cells.Item(3, 3).Value = "OverFlow:InboundPallets"
$ws.Cells.Item(3, 4).Value = "OverFlow:InboundPallets"
$ws.Range("E3:I3").ClearContents()

# --- Row 4: now the "UNLD 20 FT FLR CNT" charge row (was row 2's content) ---
$ws.Cells.Item(4, 2).Value = "UNLD 20 FT FLR CNT"
$ws.Cells.Item(4, 3).Value = "Unload 20 ft floor loaded"
$ws.Cells.Item(4, 5).Value = "WarehouseReceiptView"
$ws.Cells.Item(4, 6).Value = "TransportMethod"
$ws.Cells.Item(4, 7).Value = "20 FLR"
$ws.Cells.Item(4, 8).Value = 1
$ws.Cells.Item(4, 9).Value = 1

# --- Row 5 (new): "UNLD 40 FT FLR CNT" charge row (was row 3's content) ---
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "UNLD 40 FT FLR CNT"
$ws.Cells.Item(5, 3).Value = "Unload 40 ft floor loaded"
$ws.Cells.Item(5, 5).Value = "WarehouseReceiptView"
$ws.Cells.Item(5, 6).Value = "TransportMethod"
$ws.Cells.Item(5, 7).Value = "40 FLR"
$ws.Cells.Item(5, 8).Value = 1
$ws.Cells.Item(5, 9).Value = 1

# --- Row 6 (new): "UNLD 45 FT FLR CNT" charge row (was row 4's content) ---
$ws.Range("A4").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "UNLD 45 FT FLR CNT"
$ws.Cells.Item(6, 3).Value = "Unload 45 ft floor loaded"
$ws.Cells.Item(6, 5).Value = "WarehouseReceiptView"
$ws.Cells.Item(6, 6).Value = "TransportMethod"
$ws.Cells.Item(6, 7).Value = "45 FLR"
$ws.Cells.Item(6, 8).Value = 1
$ws.Cells.Item(6, 9).Value = 1

# --- Row 7 (new): "Vehicle Inspection-Inbound" charge row ---
$ws.Range("A4").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "Vehicle Inspection-Inbound"
$ws.Cells.Item(7, 3).Value = "Vehicle Inspection-Inbound"

$excel.CutCopyMode = 0
